$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tests")

# Update CRS point feature-class names used for the "get_crs_name" tests
# (renamed from *_LL_point to *_point as part of the CRS class refactor).
$ws.Range("F3").Value = "GDA_point"
$ws.Range("F4").Value = "WGS84_point"

# Move the active selection to reflect where the author was last working.
$ws.Range("F8").Select()
